$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column L (year 2022) mirroring the layout/styling of column K.
# Copy each source cell's formatting (and value, later overwritten) from the
# matching K-column cell so the new L-column cells pick up the same visual
# style (font, border, number format, alignment) used throughout the table.

$ws.Range("K4").Copy($ws.Range("L4"))
$ws.Range("L4").Value = 2022

$ws.Range("K5").Copy($ws.Range("L5"))
$ws.Range("L5").Value = 8800.6

$ws.Range("K6").Copy($ws.Range("L6"))

$ws.Range("K7").Copy($ws.Range("L7"))
$ws.Range("L7").Formula = "=L5-L8"

$ws.Range("K8").Copy($ws.Range("L8"))
$ws.Range("L8").Value = 258.4

$ws.Range("K9").Copy($ws.Range("L9"))

$ws.Range("K10").Copy($ws.Range("L10"))
$ws.Range("L10").Value = 683.8

$ws.Range("K11").Copy($ws.Range("L11"))
$ws.Range("L11").Value = 1101.8

$ws.Range("K12").Copy($ws.Range("L12"))
$ws.Range("L12").Value = 714.9

$ws.Range("K13").Copy($ws.Range("L13"))
$ws.Range("L13").Value = 757.9

$ws.Range("K14").Copy($ws.Range("L14"))
$ws.Range("L14").Value = 1383.3

$ws.Range("K15").Copy($ws.Range("L15"))
$ws.Range("L15").Value = 1023.7

$ws.Range("K16").Copy($ws.Range("L16"))
$ws.Range("L16").Value = 2929.3

$ws.Range("K17").Copy($ws.Range("L17"))
$ws.Range("L17").Value = 148.9

$ws.Range("K18").Copy($ws.Range("L18"))
$ws.Range("L18").Value = 57

# Update the selected cell to match the author's saved cursor position.
[void]$ws.Range("M4").Select()
